$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "6.38") are not
# auto-converted to numbers; values will match the original inlineStr text cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.761.32"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "3.009.54"
$ws.Range("E3").Value = "  -3.07%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "555.59"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "152.66"
$ws.Range("E6").Value = "  -5.02%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("D9").Value = "3.012.36"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").Value = "6.38"
$ws.Range("E11").Value = "  -4.55%  "
$ws.Range("D12").Value = "0.365"
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("D13").Value = "3.530.95"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("E14").Value = "  -3.23%  "
$ws.Range("D15").Value = "62.818.37"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "23.93"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("D17").Value = "3.011.31"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "396.32"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "11.87"
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("D22").Value = "6.62"
$ws.Range("E22").Value = "  -4.87%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "65.01"
$ws.Range("E24").Value = "  -2.85%  "
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("E26").Value = "  -5.36%  "
$ws.Range("D27").Value = "0.0₃0967"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").Value = "8.62"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "1.76"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").Value = "20.46"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").Value = "160.37"
$ws.Range("E33").Value = "  +5.71%  "
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("D38").Value = "1.58"
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("D39").Value = "2.491.39"
$ws.Range("E39").Value = "  -7.90%  "
$ws.Range("D40").Value = "37.54"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").Value = "3.91"
$ws.Range("E41").Value = "  -2.98%  "
$ws.Range("D42").Value = "22.46"
$ws.Range("E42").Value = "  -2.62%  "
$ws.Range("D43").Value = "0.664"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("D44").Value = "0.0594"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "0.0247"
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("D47").Value = "5.01"
$ws.Range("E47").Value = "  -7.74%  "
$ws.Range("D48").Value = "19.79"
$ws.Range("E48").Value = "  -4.05%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.0945"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "10.48"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "262.99"
$ws.Range("E51").Value = "  -5.38%  "

# Restore default style on column D (removes the temporary text-format style,
# keeping cells styled the same as before the edit).
$ws.Range("D2:D51").Style = "Normal"
